$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep their literal text representation
# (e.g. trailing zeros, thousand-dot formatting) instead of Excel auto-
# converting numeric-looking strings into real numbers when assigned.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.083.50"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").Value = "1.909.61"
$ws.Range("E3").Value = "  +0.81%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").Value = "319.83"
$ws.Range("E5").Value = "  -0.78%  "
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").Value = "0.5035"
$ws.Range("E7").Value = "  -2.34%  "
$ws.Range("D8").Value = "0.4074"
$ws.Range("E8").Value = "  +1.80%  "
$ws.Range("D9").Value = "0.08325"
$ws.Range("E9").Value = "  -1.02%  "
$ws.Range("D10").Value = "42.17"
$ws.Range("E10").Value = "  -0.86%  "
$ws.Range("D11").Value = "1.102"
$ws.Range("E11").Value = "  -0.71%  "
$ws.Range("D12").Value = "23.93"
$ws.Range("E12").Value = "  +2.94%  "
$ws.Range("D13").Value = "1.910.23"
$ws.Range("E13").Value = "  +1.05%  "
$ws.Range("D14").Value = "6.384"
$ws.Range("E14").Value = "  -0.28%  "
$ws.Range("D15").Value = "7.211"
$ws.Range("E15").Value = "  -0.90%  "
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("D17").Value = "92.20"
$ws.Range("E17").Value = "  -2.14%  "
$ws.Range("D18").Value = "0.00001095"
$ws.Range("E18").Value = "  -0.88%  "
$ws.Range("D19").Value = "0.06487"
$ws.Range("E19").Value = "  -0.43%  "
$ws.Range("D20").Value = "18.26"
$ws.Range("E20").Value = "  +0.48%  "
$ws.Range("D21").Value = "1.001"
$ws.Range("E22").Value = "  +0.34%  "
$ws.Range("D23").Value = "30.074.98"
$ws.Range("E23").Value = "  -0.61%  "
$ws.Range("D24").Value = "11.30"
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").Value = "2.191"
$ws.Range("E25").Value = "  -1.58%  "
$ws.Range("D26").Value = "2.130.81"
$ws.Range("E26").Value = "  +1.09%  "
$ws.Range("D27").Value = "21.71"
$ws.Range("E27").Value = "  -0.31%  "
$ws.Range("D28").Value = "163.02"
$ws.Range("E28").Value = "  +0.58%  "
$ws.Range("D29").Value = "2.289"
$ws.Range("E29").Value = "  -1.15%  "
$ws.Range("D30").Value = "128.61"
$ws.Range("E30").Value = "  -0.74%  "
$ws.Range("D31").Value = "1.153"
$ws.Range("E31").Value = "  +7.49%  "
$ws.Range("D32").Value = "0.1038"
$ws.Range("E32").Value = "  -0.79%  "
$ws.Range("D33").Value = "5.957"
$ws.Range("E33").Value = "  -1.33%  "
$ws.Range("D34").Value = "3.727"
$ws.Range("E34").Value = "  -4.19%  "
$ws.Range("D35").Value = "0.02449"
$ws.Range("E35").Value = "  -0.94%  "
$ws.Range("D36").Value = "5.371"
$ws.Range("E36").Value = "  +1.38%  "
$ws.Range("D37").Value = "0.06361"
$ws.Range("E37").Value = "  -1.48%  "
$ws.Range("D38").Value = "0.2141"
$ws.Range("E38").Value = "  -2.09%  "
$ws.Range("D39").Value = "0.6545"
$ws.Range("E39").Value = "  +1.76%  "
$ws.Range("D40").Value = "1.187"
$ws.Range("E40").Value = "  -0.96%  "
$ws.Range("D41").Value = "8.627"
$ws.Range("E41").Value = "  -0.80%  "
$ws.Range("D42").Value = "11.35"
$ws.Range("E42").Value = "  -2.93%  "
$ws.Range("D43").Value = "1.211"
$ws.Range("E43").Value = "  -1.27%  "
$ws.Range("D46").Value = "0.6059"
$ws.Range("E46").Value = "  +0.32%  "
$ws.Range("D47").Value = "3.622"
$ws.Range("E47").Value = "  -1.12%  "
$ws.Range("E48").Value = "  -1.74%  "
$ws.Range("D49").Value = "121.31"
$ws.Range("E49").Value = "  -2.16%  "
$ws.Range("D50").Value = "78.83"
$ws.Range("E50").Value = "  +0.52%  "
$ws.Range("D51").Value = "1.141"
$ws.Range("E51").Value = "  +0.26%  "

# Rows 44/45 swap places: EnergySwap <-> NEARProtocol
$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").Value = "2.214"
$ws.Range("E44").Value = "  +8.65%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "13.42"
$ws.Range("E45").Value = "  +2.68%  "

Write-Host "Applied cryptos update"
